$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 185 (pushes the existing rows 185-221 down to 186-222,
# carrying formatting such as the date style on column D along with them).
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new weekly record.
$ws.Cells.Item(185, 1).Value = 4
$ws.Cells.Item(185, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(185, 3).Value = "Los Lagos"
$ws.Cells.Item(185, 4).Value = 44641
$ws.Cells.Item(185, 5).Value = 10
$ws.Cells.Item(185, 6).Value = 100112017
$ws.Cells.Item(185, 7).Value = "Apio"
$ws.Cells.Item(185, 8).Value = "Americana (o)"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 25
$ws.Cells.Item(185, 11).Value = 13000
$ws.Cells.Item(185, 12).Value = 13000
$ws.Cells.Item(185, 13).Value = 13000
$ws.Cells.Item(185, 14).Value = "`$/docena de matas"
$ws.Cells.Item(185, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(185, 16).Value = 2167
$ws.Cells.Item(185, 17).Value = 6
$ws.Cells.Item(185, 18).Value = "Hortaliza"
